$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5,6,7 (the old FAPs-sender->ECs-target and MuSCs-sender->ECs-target rows,
# and the old FAPs->FAPs row that shifts up) - we will instead rewrite rows 2-4 entirely
# and delete rows 5-7.

# Row 2: ECs -> FAPs (was ECs -> ECs)
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4055383333333333
$ws.Range("N2").Value = 1.216615
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 33.65468583993833
$ws.Range("R2").Value = 302.892172559445
$ws.Range("S2").Value = 0.4489504115427952
$ws.Range("T2").Value = 0.4489504115427952

# Row 3: FAPs -> FAPs (was ECs -> FAPs)
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 63.14058933333333
$ws.Range("H3").Value = 189.421768
$ws.Range("I3").Value = 0.3415807409566563
$ws.Range("J3").Value = 0.3415807409566563
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 25.60592936392444
$ws.Range("R3").Value = 230.45336427532
$ws.Range("S3").Value = 0.3415807409566563
$ws.Range("T3").Value = 0.3415807409566563

# Row 4: MuSCs -> FAPs (was FAPs -> ECs)
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 38.719942
$ws.Range("H4").Value = 116.159826
$ws.Range("I4").Value = 0.2094688475005485
$ws.Range("J4").Value = 0.2094688475005485
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4055383333333333
$ws.Range("N4").Value = 1.216615
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 15.70242074544333
$ws.Range("R4").Value = 141.32178670899
$ws.Range("S4").Value = 0.2094688475005485
$ws.Range("T4").Value = 0.2094688475005485

# Delete old rows 5, 6, 7 (FAPs->FAPs duplicate, MuSCs->ECs, MuSCs->FAPs) which are no longer needed
$ws.Range("A5:T7").Delete()
